$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 106; existing rows 106-122 shift down to 107-123
$ws.Rows.Item(106).Insert()

# Populate the new row 106 with the new record
$ws.Cells.Item(106, 1).Value = 10
$ws.Cells.Item(106, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(106, 3).Value = "La Araucanía"
$ws.Cells.Item(106, 4).Value = 44617
$ws.Cells.Item(106, 4).NumberFormat = $ws.Cells.Item(107, 4).NumberFormat
$ws.Cells.Item(106, 5).Value = 9
$ws.Cells.Item(106, 6).Value = "Fruta"
$ws.Cells.Item(106, 7).Value = 100104
$ws.Cells.Item(106, 8).Value = "Frutos de pepita"
$ws.Cells.Item(106, 9).Value = 100104003
$ws.Cells.Item(106, 10).Value = "Membrillo"
$ws.Cells.Item(106, 11).Value = "Champion"
$ws.Cells.Item(106, 12).Value = "Primera"
$ws.Cells.Item(106, 13).Value = 20
$ws.Cells.Item(106, 14).Value = 16000
$ws.Cells.Item(106, 15).Value = 16000
$ws.Cells.Item(106, 16).Value = 16000
$ws.Cells.Item(106, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(106, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(106, 19).Value = 889
$ws.Cells.Item(106, 20).Value = 18
